$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Neo4j Cypher query text stored in the "FilesTab" query cell (B4):
# the "File Type" and "Breed" columns were dropped from the RETURN clause.
$ws.Range("B4").Value = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['GLIOMA01'] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# The author's saved selection moved from C4 to B4.
$ws.Range("B4").Select() | Out-Null
